# PlayerPerformance_4393.xlsx update
# - Insert a new "Player Info" sheet at the front.
# - Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling",
#   replacing the full howstat URL values with the bare numeric match code.
# - Drop the empty INNING_NUMBER (column B) cells on "ODI Batting".
# - Append a new "ODI Batting Extra" sheet at the end with additional
#   per-match batting detail.

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# 1. "Player Info" sheet, inserted before "ODI Batting".
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$piHeader = $playerInfo.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.HorizontalAlignment = -4108
$piHeader.VerticalAlignment = -4160
$piHeader.Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4393"
$playerInfo.Range("B2").Value = "Shardul Narendra Thakur"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, url -> bare code,
#    drop empty INNING_NUMBER cells.
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingRows = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $battingRows; $r++) {
    $linkCell = $battingSheet.Cells.Item($r, 4)
    $linkText = $linkCell.Text
    if ($linkText -match "MatchCode=(\d+)") {
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $matches[1]
    }

    $inningCell = $battingSheet.Cells.Item($r, 2)
    if ($inningCell.Text -eq "") {
        $inningCell.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, url -> bare code.
# ---------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingRows = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingRows; $r++) {
    $linkCell = $bowlingSheet.Cells.Item($r, 2)
    $linkText = $linkCell.Text
    if ($linkText -match "MatchCode=(\d+)") {
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $matches[1]
    }
}

# ---------------------------------------------------------------------------
# 4. "ODI Batting Extra", appended after "ODI Bowling".
# ---------------------------------------------------------------------------
$extra = $wb.Worksheets.Add($null, $bowlingSheet)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

$extraHeader = $extra.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160
$extraHeader.Borders.LineStyle = 1

function Set-ExtraRow($row, $matchCode, $battingPosition, $num4, $num6, $percent, $manOfMatch) {
    $codeCell = $extra.Cells.Item($row, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $matchCode

    if ($null -ne $battingPosition) {
        $extra.Cells.Item($row, 2).Value = $battingPosition
    }

    if ($null -ne $num4) {
        $c4 = $extra.Cells.Item($row, 3)
        $c4.NumberFormat = "@"
        $c4.Value = $num4
    }

    if ($null -ne $num6) {
        $c6 = $extra.Cells.Item($row, 4)
        $c6.NumberFormat = "@"
        $c6.Value = $num6
    }

    if ($null -ne $percent) {
        $cp = $extra.Cells.Item($row, 5)
        $cp.NumberFormat = "@"
        $cp.Value = $percent
    }

    $extra.Cells.Item($row, 6).Value = $manOfMatch
}

Set-ExtraRow 2  "4524" $null $null $null $null "NO"
Set-ExtraRow 3  "4526" 7    "3"  "1"  "13.94%" "NO"
Set-ExtraRow 4  "4533" $null $null $null $null "NO"
Set-ExtraRow 5  "4535" 8    "1"  "0"  "3.38%"  "NO"
Set-ExtraRow 6  "4621" 8    "1"  "0"  "2.27%"  "NO"
Set-ExtraRow 7  "4623" $null $null $null $null "NO"
Set-ExtraRow 8  "4624" 8    $null $null $null  "NO"
Set-ExtraRow 9  "4640" 8    $null $null $null  "NO"
Set-ExtraRow 10 "4643" $null $null $null $null "NO"
Set-ExtraRow 11 "4656" $null $null $null $null "NO"
Set-ExtraRow 12 "4657" 6    $null $null $null  "NO"
Set-ExtraRow 13 "4658" 6    $null $null $null  "NO"
Set-ExtraRow 14 "4669" 8    "0"  "0"  "0.33%"  "NO"
Set-ExtraRow 15 "4679" 8    "0"  "0"  "1.08%"  "NO"
Set-ExtraRow 16 "4682" 7    "0"  "0"  "2.63%"  "NO"
Set-ExtraRow 17 "4685" $null $null $null $null "NO"
Set-ExtraRow 18 "4692" $null $null $null $null "NO"
Set-ExtraRow 19 "4695" 8    $null $null $null  "NO"
Set-ExtraRow 20 "4697" 8    "3"  "1"  "6.49%"  "YES"
Set-ExtraRow 21 "4725" 8    $null $null $null  "NO"

# ---------------------------------------------------------------------------
# 5. Keep the first sheet active (matches activeTab="0").
# ---------------------------------------------------------------------------
$playerInfo.Activate()
